# The two sightings previously stored in rows 3 and 4 were re-ordered:
# everything that used to be on row 4 now belongs on row 3 and vice
# versa. Columns C, P, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY
# already hold identical values on both rows, so only the columns below
# actually need to move.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteAll = -4104

foreach ($col in @("A", "B", "D", "E", "F", "G", "H", "I", "Q", "R", "S", "Z", "AB")) {
    $addr3 = $col + "3"
    $addr4 = $col + "4"
    $scratch = $col + "100"

    # Use Copy/PasteSpecial (instead of reading .Value() and writing it
    # back) so a cell's original storage type survives the swap - e.g.
    # column I holds the text "35" on one row, which a Value round trip
    # would otherwise silently reinterpret as a number. Each
    # destination is cleared first because pasting a blank source cell
    # over a non-blank one otherwise leaves the old content in place.
    $ws.Range($scratch).ClearContents()
    $ws.Range($addr3).Copy()
    $ws.Range($scratch).PasteSpecial($xlPasteAll)

    $ws.Range($addr3).ClearContents()
    $ws.Range($addr4).Copy()
    $ws.Range($addr3).PasteSpecial($xlPasteAll)

    $ws.Range($addr4).ClearContents()
    $ws.Range($scratch).Copy()
    $ws.Range($addr4).PasteSpecial($xlPasteAll)
}

$ws.Rows.Item(100).Delete()
